# Update the "want to go" counts (column F) that changed between scrapes
# on the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 527
$wsExhibit.Range("F11").Value = 1387
$wsExhibit.Range("F24").Value = 2
$wsExhibit.Range("F26").Value = 87
$wsExhibit.Range("F33").Value = 122
$wsExhibit.Range("F36").Value = 321

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 527
$wsAll.Range("F12").Value = 1389
$wsAll.Range("F25").Value = 2
$wsAll.Range("F27").Value = 87
$wsAll.Range("F34").Value = 122
$wsAll.Range("F37").Value = 321
